# Apply updated "dSF" (column F) values as part of a data repull / push.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    4  = 2
    12 = -5
    14 = -1
    28 = -1
    29 = -2
    39 = 0
    41 = 2
    44 = -2
    45 = 2
    50 = 1
    55 = -2
    56 = -1
    58 = -2
    62 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
